# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" fund-holdings sheet between "2021-Q2" and "总计".
# - Refresh "总计" with a new row summarising the 2022-Q1 data (pushing the
#   existing 2021-Q2 row down).

$wb = $excel.ActiveWorkbook

$sheetQ2 = $wb.Worksheets.Item("2021-Q2")
$oldTotal = $wb.Worksheets.Item("总计")

# Duplicate the existing "总计" sheet *before* touching it, placing the copy at the
# end of the workbook - this becomes the refreshed "总计" sheet, and keeps all of
# the original sheet's structural bits (sheetPr/pageMargins/etc.) intact.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldTotal.Copy($null, $lastSheet)
$newTotal = $wb.Worksheets.Item($wb.Worksheets.Count)
$newTotal.Name = "总计-new"

# The original "总计" sheet becomes the new "2022-Q1" sheet in place, so it keeps its
# position (right after "2021-Q2") and its existing formatting.
$newQ1 = $oldTotal
$newQ1.Name = "2022-Q1"

# Copy the header-row formatting into the new columns (E:H) that this sheet didn't
# have before, and copy the A2 "index" cell formatting down onto the new A3 row.
$newQ1.Range("D1").Copy()
$newQ1.Range("E1:H1").PasteSpecial(-4122)
$newQ1.Range("A2").Copy()
$newQ1.Range("A3").PasteSpecial(-4122)

# --- Header row ---
$newQ1.Range("B1").Value = "基金代码"
$newQ1.Range("C1").Value = "基金名称"
$newQ1.Range("D1").Value = "基金规模"
$newQ1.Range("E1").Value = "股票总仓位"
$newQ1.Range("F1").Value = "仓位占比"
$newQ1.Range("G1").Value = "持有市值(亿元)"
$newQ1.Range("H1").Value = "仓位排名"

# --- Data rows (numeric-looking values stored as text, matching the source data) ---
# NumberFormat only honours the first area of a multi-area Range, so set it per block.
$newQ1.Range("B2:B3").NumberFormat = "@"
$newQ1.Range("D2:G3").NumberFormat = "@"

$newQ1.Range("A2").Value = 0
$newQ1.Range("B2").Value = "003720"
$newQ1.Range("C2").Value = "易方达标普生物科技指数（QDII-LOF）美元"
$newQ1.Range("D2").Value = "2.11"
$newQ1.Range("E2").Value = "94.00"
$newQ1.Range("F2").Value = "1.02"
$newQ1.Range("G2").Value = "0.0215"
$newQ1.Range("H2").Value = 3

$newQ1.Range("A3").Value = 1
$newQ1.Range("B3").Value = "161127"
$newQ1.Range("C3").Value = "易方达标普生物科技指数（QDII-LOF）人民币"
$newQ1.Range("D3").Value = "2.11"
$newQ1.Range("E3").Value = "94.00"
$newQ1.Range("F3").Value = "1.02"
$newQ1.Range("G3").Value = "0.0215"
$newQ1.Range("H3").Value = 3

$newQ1.Range("B2:B3").ClearFormats()
$newQ1.Range("D2:G3").ClearFormats()

# --- Refresh the (copied) "总计" sheet: insert the 2022-Q1 row, shift 2021-Q2 down ---
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q2"
$newTotal.Range("C3").Value = 2
$newTotal.Range("D3").Value = 0.03

$newTotal.Range("A2").Copy()
$newTotal.Range("A3").PasteSpecial(-4122)
$newTotal.Range("A3").Value = 1

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("D2").Value = 0.04

$newTotal.Name = "总计"

# Keep the original active sheet/tab selection ("2021-Q2") as it was before the edit.
$sheetQ2.Activate()
$sheetQ2.Range("A1").Select() | Out-Null
